# Apply TPM-update edit to Fgf5-Fgfr4 LR-pairs sheet
# Adds a new "Inflammatory-Mac" cluster (target + sending) and refreshes
# expression-weight statistics with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 0.7825479339666588
$ws.Range("M2").Value = 0.903185
$ws.Range("N2").Value = 2.709555
$ws.Range("O2").Value = 0.03154869388788047
$ws.Range("P2").Value = 0.03154869388788046
$ws.Range("Q2").Value = 0.101654374935
$ws.Range("R2").Value = 0.914889374415
$ws.Range("S2").Value = 0.02468836522130742
$ws.Range("T2").Value = 0.02468836522130741

# Row 3
$ws.Range("J3").Value = 0.7825479339666588
$ws.Range("O3").Value = 0.003855347953955327
$ws.Range("P3").Value = 0.003855347953955326
$ws.Range("S3").Value = 0.003016994576090327
$ws.Range("T3").Value = 0.003016994576090326

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("J4").Value = 0.7825479339666588
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03076233333333333
$ws.Range("N4").Value = 0.09228699999999999
$ws.Range("O4").Value = 0.001074543352259254
$ws.Range("P4").Value = 0.001074543352259254
$ws.Range("Q4").Value = 0.003462331379
$ws.Range("R4").Value = 0.03116098241099999
$ws.Range("S4").Value = 0.0008408816802680873
$ws.Range("T4").Value = 0.0008408816802680871

# Row 5
$ws.Range("D5").Value = "MuSCs"
$ws.Range("J5").Value = 0.7825479339666588
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 27.52907633333334
$ws.Range("N5").Value = 82.58722900000001
$ws.Range("O5").Value = 0.9616041035407232
$ws.Range("P5").Value = 0.9616041035407231
$ws.Range("Q5").Value = 3.098425070393
$ws.Range("R5").Value = 27.885825633537
$ws.Range("S5").Value = 0.7525013045196541
$ws.Range("T5").Value = 0.7525013045196539

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("G6").Value = 0.112551
$ws.Range("H6").Value = 0.337653
$ws.Range("I6").Value = 0.7825479339666589
$ws.Range("J6").Value = 0.7825479339666588
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05488933333333334
$ws.Range("N6").Value = 0.164668
$ws.Range("O6").Value = 0.001917311265181737
$ws.Range("P6").Value = 0.001917311265181736
$ws.Range("Q6").Value = 0.006177849356
$ws.Range("R6").Value = 0.055600644204
$ws.Range("S6").Value = 0.001500387969338969
$ws.Range("T6").Value = 0.001500387969338969

# Row 7
$ws.Range("D7").Value = "ECs"
$ws.Range("G7").Value = 0.03127533333333334
$ws.Range("H7").Value = 0.09382600000000001
$ws.Range("I7").Value = 0.2174520660333412
$ws.Range("J7").Value = 0.2174520660333412
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.903185
$ws.Range("N7").Value = 2.709555
$ws.Range("O7").Value = 0.03154869388788047
$ws.Range("P7").Value = 0.03154869388788046
$ws.Range("Q7").Value = 0.02824741193666667
$ws.Range("R7").Value = 0.25422670743
$ws.Range("S7").Value = 0.00686032866657305
$ws.Range("T7").Value = 0.006860328666573048

# Row 8
$ws.Range("D8").Value = "FAPs"
$ws.Range("G8").Value = 0.03127533333333334
$ws.Range("H8").Value = 0.09382600000000001
$ws.Range("I8").Value = 0.2174520660333412
$ws.Range("J8").Value = 0.2174520660333412
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.110372
$ws.Range("N8").Value = 0.331116
$ws.Range("O8").Value = 0.003855347953955327
$ws.Range("P8").Value = 0.003855347953955326
$ws.Range("Q8").Value = 0.003451921090666667
$ws.Range("R8").Value = 0.031067289816
$ws.Range("S8").Value = 0.0008383533778650005
$ws.Range("T8").Value = 0.0008383533778650004

# Row 9
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("G9").Value = 0.03127533333333334
$ws.Range("H9").Value = 0.09382600000000001
$ws.Range("I9").Value = 0.2174520660333412
$ws.Range("J9").Value = 0.2174520660333412
$ws.Range("M9").Value = 0.03076233333333333
$ws.Range("N9").Value = 0.09228699999999999
$ws.Range("O9").Value = 0.001074543352259254
$ws.Range("P9").Value = 0.001074543352259254
$ws.Range("Q9").Value = 0.0009621022291111111
$ws.Range("R9").Value = 0.008658920061999999
$ws.Range("S9").Value = 0.0002336616719911672
$ws.Range("T9").Value = 0.0002336616719911672

# Row 10 (new)
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf5"
$ws.Range("C10").Value = "Fgfr4"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.03127533333333334
$ws.Range("H10").Value = 0.09382600000000001
$ws.Range("I10").Value = 0.2174520660333412
$ws.Range("J10").Value = 0.2174520660333412
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 27.52907633333334
$ws.Range("N10").Value = 82.58722900000001
$ws.Range("O10").Value = 0.9616041035407232
$ws.Range("P10").Value = 0.9616041035407231
$ws.Range("Q10").Value = 0.860981038683778
$ws.Range("R10").Value = 7.748829348154001
$ws.Range("S10").Value = 0.2091027990210692
$ws.Range("T10").Value = 0.2091027990210692

# Row 11 (new)
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Fgf5"
$ws.Range("C11").Value = "Fgfr4"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.03127533333333334
$ws.Range("H11").Value = 0.09382600000000001
$ws.Range("I11").Value = 0.2174520660333412
$ws.Range("J11").Value = 0.2174520660333412
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.05488933333333334
$ws.Range("N11").Value = 0.164668
$ws.Range("O11").Value = 0.001917311265181737
$ws.Range("P11").Value = 0.001917311265181736
$ws.Range("Q11").Value = 0.001716682196444445
$ws.Range("R11").Value = 0.015450139768
$ws.Range("S11").Value = 0.0004169232958427679
$ws.Range("T11").Value = 0.0004169232958427679
